$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells with the new "c" value
$ws.Range("A9").Value = "c"
$ws.Range("B10").Value = "c"

# Add new cells (rows 11-15) with the new "c" value
$ws.Range("C11").Value = "c"
$ws.Range("D11").Value = "c"
$ws.Range("B12").Value = "c"
$ws.Range("D12").Value = "c"
$ws.Range("B13").Value = "c"
$ws.Range("B14").Value = "c"
$ws.Range("D14").Value = "c"
$ws.Range("B15").Value = "c"
$ws.Range("D15").Value = "c"

# Update the active selection
$ws.Range("A9").Select()
